$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cadastro")

# Update the username value (shared string "john1115" -> "john1116")
# This string appears in both B2 and B14 on sheet "Cadastro"
$ws.Range("B2").Value = "john1116"
$ws.Range("B14").Value = "john1116"

# Move the active selection on sheet "Cadastro" to B14
$ws.Range("B14").Select()
